$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: in the "Interpreter ... design pattern" bullet, turn
#   " to translate objects pulled from the database to objects used
#    in our backend code. "
# into an (underlined) colon immediately after "design pattern",
# followed by a capitalised, non-underlined continuation:
#   ": To translate objects pulled from the database to objects used
#    in our backend code. "
# -----------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute(
    " to translate objects pulled from the database to objects used in our backend code. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $find.Start

    # Replace the located text in-place (keeps the run's original,
    # non-underlined formatting) with the capitalised continuation.
    $find.Text = " To translate objects pulled from the database to objects used in our backend code. "

    # Insert the new underlined colon immediately before it.
    $colonSpot = $d.Range($start, $start)
    $colonSpot.InsertBefore(":")

    # Underline only the colon we just added.
    $colonRange = $d.Range($start, $start + 1)
    $colonRange.Font.Underline = 1
}

# -----------------------------------------------------------------
# Change 2: in the (previously empty) last bullet of the
# "Design Patterns" list, add a new description of the Singleton
# pattern before the existing _GoBack bookmark.
# -----------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$labelStart = $bm.Range.Start

$segments = @(
    "Singleton creational design pattern",
    ": ",
    "We made most of the methods in the classes (",
    "DatabaseTranslator",
    " and ",
    "DatabaseManagerImp",
    ") static so that their methods can be used throughout the program instead of having to make an instance every time."
)

foreach ($segment in $segments) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Range.InsertBefore($segment)
}

# Underline just the "Singleton creational design pattern" label.
$labelRange = $d.Range($labelStart, $labelStart + "Singleton creational design pattern".Length)
$labelRange.Font.Underline = 1
